# "format float in psgq" - update the ticker (CRWD -> DDOG) and its
# associated metrics on both report sheets, and reformat the
# Price/Sales/Growth float series to 2-decimal precision.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # IncomeReport
$ws2 = $wb.Worksheets.Item(2)   # EarningsReport

# --- IncomeReport (sheet1) -------------------------------------------------

# Ticker (shared with EarningsReport!A2 through the shared-string table)
$ws1.Range("A2").Value = "DDOG"

# Revenue / Rev growth / Earnings / EPS growth / FCF / FCF growth series
$ws1.Range("B2").Value = "436.53M -> 469.4M -> 481.71M -> 509.46M -> 547.54M"
$ws1.Range("C2").Value = 25
$ws1.Range("D2").Value = "-0.08 -> -0.09 -> -0.08 -> -0.01 -> 0.06"
$ws1.Range("E2").Value = 175
$ws1.Range("F2").Value = "73.91M -> 104.39M -> 125.05M -> 150.82M -> 146.67M"
$ws1.Range("G2").Value = 98

# Price/Earnings/Growth, Price/Sales/Growth (now 2dp formatted), Price/Sales
$ws1.Range("H2").Value = "1.31 <- 2.05 <- 1.69 <- 1.12 <- 1.49"
$ws1.Range("I2").Value = "0.61 <- 0.70 <- 0.55 <- 0.60 <- 0.85"
$ws1.Range("J2").Value = "15.30 <- 17.38 <- 13.68 <- 15.07 <- 21.34"

# Column width tweaks that came with the re-flowed text
$ws1.Columns.Item(2).ColumnWidth = 24.8333333333333
$ws1.Columns.Item(4).ColumnWidth = 19.8333333333333
$ws1.Columns.Item(6).ColumnWidth = 24.8333333333333
$ws1.Columns.Item(9).ColumnWidth = 17.8333333333333

# --- EarningsReport (sheet2) ------------------------------------------------

# Ticker stays in lock-step with IncomeReport (same underlying shared string)
$ws2.Range("A2").Value = "DDOG"

# EPS Estimate / EPS Actual must remain text, not be coerced to numbers
$ws2.Range("B2").Value = "'0.34"
$ws2.Range("C2").Value = "'0.45"

$ws2.Range("D2").Value = 32
$ws2.Range("E2").Value = "633, 1056"
$ws2.Range("F2").Value = "21, 25"
$ws2.Range("G2").Value = 22
